# "fixed export and fixing maps"
#
# The sheet had a subtitle row ("(by census results)") that no longer
# applies, and two of the three year columns (1989, 2002) are dropped,
# leaving only the most recent (2014) figures. The sheet is also given
# its real Georgian name instead of the generic "1" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unwanted subtitle row (row 2) - rows below shift up.
$ws.Rows("2:2").Delete()

# Keep only the 2014 figures - drop the 1989 and 2002 columns (B and C);
# the former "2014" column (D) slides left into column B.
$ws.Columns("B:C").Delete()

# Give the sheet/tab its proper name.
$ws.Name = "ახალქალაქი"

# Leave the cursor on A2, matching the saved selection.
$ws.Range("A2").Select()
